# Weekly update: a new daily price record was added for
# "Vega Modelo de Temuco - Bruselas (repollito)".
# The new record is inserted as row 17 (rows formerly at 17..121 move
# down to 18..122), matching the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 17; this shifts the existing
# rows 17-121 down to 18-122 (and extends the used range to R122).
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(17, 1).Value  = 10
$ws.Cells.Item(17, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value  = "La Araucanía"
$ws.Cells.Item(17, 4).Value  = 44802
$ws.Cells.Item(17, 5).Value  = 9
$ws.Cells.Item(17, 6).Value  = 100112035
$ws.Cells.Item(17, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(17, 8).Value  = "Sin especificar"
$ws.Cells.Item(17, 9).Value  = "Primera"
$ws.Cells.Item(17, 10).Value = 90
$ws.Cells.Item(17, 11).Value = 24000
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24389
$ws.Cells.Item(17, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 16).Value = 2439
$ws.Cells.Item(17, 17).Value = 10
$ws.Cells.Item(17, 18).Value = "Hortaliza"
